# Fill all 21 ActivityLogs columns including LoginHour, RiskLabel, Channel and ML fields

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ActivityLogs")

$headers = @(
    "LogID",
    "AccountID",
    "Timestamp",
    "TransactionType",
    "Description",
    "SessionID",
    "TransactionAmount",
    "SessionDuration",
    "LoginHour",
    "FailedLoginCount",
    "NewDeviceLogin",
    "PasswordChanged",
    "Channel",
    "PagesVisited",
    "ClickRate",
    "RapidTransactions",
    "BeneficiaryAdded",
    "LargeTransaction",
    "DeviceTrustScore",
    "CyberRiskScore",
    "RiskLabel"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}
